$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column D values from "N" to "Y" for rows 2 through 22
$ws.Range("D2:D22").Value = "Y"

# Update the selection to match the new range
$ws.Range("D2:D22").Select()
